$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Topic labels for rows 261-330 (column D), in row order.
$topics = @("مساعدة","قروض","قروض","مساعدة","مساعدة","مساعدة","مساعدة","مساعدة","مساعدة","مساعدة","مساعدة","بطاقات","قروض","قروض","حسابات","حسابات","بطاقات","بطاقات","حسابات","حسابات","حسابات","حسابات","حسابات","حسابات","ودائع","حسابات","شهادات","شهادات","قروض","قروض","حسابات","حسابات","حسابات","حسابات","حسابات","مساعدة","بطاقات","بطاقات","بطاقات","قروض","قروض","ودائع","حسابات","بطاقات","قروض","بطاقات","بطاقات","شهادات","قروض","قروض","عقاري","مساعدة","مساعدة","قروض","عقاري","عقاري","حسابات","حسابات","قروض","قروض","قروض","قروض","ودائع","ودائع","ودائع","ودائع","بطاقات","قروض","قروض","قروض")

$startRow = 261
for ($i = 0; $i -lt $topics.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value = $topics[$i]
}

# Restore the author's view state: scrolled down with F330 selected.
$ws.Range("F330").Select()
